# Update the "Przygotowanie środowiska" slide: split the first sentence's
# run so "Sklonuj " / "repozytorium " are separate runs, and swap the
# repository hyperlink (both its display text and its target address) for
# the new one, dropping the stray trailing-space run that used to sit after
# the link.

$p = $ppt.ActivePresentation

$oldUrl = "https://github.com/TheVosges/Article-Classifier"
$newUrl = "https://github.com/piotrbarabasz/WebArticleClassifier?tab=readme-ov-file"
$splitAfter = "Sklonuj "

# Locate the shape that holds the hyperlinked repository URL instead of
# assuming a fixed slide/shape index.
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text.IndexOf($oldUrl) -ge 0) {
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$full = $tr.Text

# --- 1) Split "Sklonuj repozytorium " into "Sklonuj " + "repozytorium " ---
$splitPos = $full.IndexOf($splitAfter) + $splitAfter.Length
$firstPart = $tr.Characters(1, $splitPos)
$firstPart.Text = $firstPart.Text

# --- 2) Swap the URL run's visible text for the new address ---
$full = $tr.Text
$urlStart = $full.IndexOf($oldUrl) + 1
$urlRange = $tr.Characters($urlStart, $oldUrl.Length)
$urlRange.Text = $newUrl

# --- 3) Point the hyperlink itself at the new address ---
$newUrlRange = $tr.Characters($urlStart, $newUrl.Length)
$hlink = $newUrlRange.ActionSettings(1).Hyperlink
$hlink.Address = $newUrl

# --- 4) Drop the leftover single-space run right after the link ---
$full = $tr.Text
$spacePos = $urlStart + $newUrl.Length
if ($spacePos -le $full.Length -and $full.Substring($spacePos - 1, 1) -eq " ") {
    $trailingSpace = $tr.Characters($spacePos, 1)
    $trailingSpace.Text = ""
}

Write-Host "Updated text: $($tr.Text)"
Write-Host "Updated hyperlink address: $($hlink.Address)"
